$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3146293670956112
$ws.Range("C2").Value = 0.003648818996706889

$ws.Range("B3").Value = 0.3004422598709933
$ws.Range("C3").Value = 0.001959607579694902

$ws.Range("B4").Value = 0.2586819355893186
$ws.Range("C4").Value = 0.003922771422438028

$ws.Range("B5").Value = 0.3060222718072834
$ws.Range("C5").Value = 0.004555794877035837

$ws.Range("B6").Value = 0.2041455830938401
$ws.Range("C6").Value = 0.002344272394948808

$ws.Range("B7").Value = 0.1122919010182405
$ws.Range("C7").Value = 0.004190936201316352

$ws.Range("B8").Value = 0.04295661193951493
$ws.Range("C8").Value = 0.00134131295413297
